# Applies the edit described by the commit: "prot date fixed, docs form pre fixed"
#
# Summary of the semantic text changes inside the body of the document:
#   1. "№ 1-1015705121816-Д" -> "№ 1-{docNumber}-Д"   (hard-coded protocol number -> merge placeholder)
#   2. "{taskNum}"           -> "{cNumber}"             (placeholder renamed)
#   3. "δQn = мінус {Qn}%; δQt = мінус {Qt}%; δQmin = мінус {Qmin}%."
#                             -> "{testsVal}."           (three placeholders collapsed into one)

$d = $word.ActiveDocument

# 1) Replace the hard-coded protocol/document number with the {docNumber} placeholder.
$d.Content.Find.Execute("1015705121816", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{docNumber}", 2)

# 2) Rename the {taskNum} placeholder to {cNumber}.
$d.Content.Find.Execute("{taskNum}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{cNumber}", 2)

# 3) Collapse the three relative-error placeholders into the single {testsVal} placeholder.
$d.Content.Find.Execute("δQn = мінус {Qn}%; δQt = мінус {Qt}%; δQmin = мінус {Qmin}%.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{testsVal}.", 2)
